# Apply the "What do you like_analysis.xlsx" update:
#  - A couple of survey responses (row 8 & row 9) were corrected to a
#    different phone / preferred president.
#  - The "Frequency" column (H) was recomputed for several rows.
#  - The header row on the "Data" sheet got a taller row height (the
#    header formatting itself - bold font, borders, centered/top
#    alignment - stays the same).
# The same corrections are mirrored on both the "Data" sheet and the
# "Sheet1" sheet (they hold duplicate copies of the same table).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Data", "Sheet1")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- corrected survey answers (row 8 and row 9) ---
    $ws.Range("F8").Value = "iPhone SE 2020"
    $ws.Range("G8").Value = "Barack Obama"

    $ws.Range("F9").Value = "iPhone 8"
    $ws.Range("G9").Value = "Kanye West"

    # --- recomputed "Frequency" values ---
    $ws.Range("H4").Value  = 12
    $ws.Range("H6").Value  = 32
    $ws.Range("H7").Value  = 32
    $ws.Range("H8").Value  = 19
    $ws.Range("H9").Value  = 25
    $ws.Range("H10").Value = 25
    $ws.Range("H11").Value = 27
    $ws.Range("H13").Value = 8
    $ws.Range("H14").Value = 16
    $ws.Range("H15").Value = 16
    $ws.Range("H17").Value = 26
    $ws.Range("H19").Value = 32
    $ws.Range("H20").Value = 16
    $ws.Range("H23").Value = 21
    $ws.Range("H24").Value = 11
    $ws.Range("H25").Value = 10
    $ws.Range("H26").Value = 15
    $ws.Range("H28").Value = 24
}

# --- header row on the "Data" sheet gets a taller row height ---
$dataWs = $wb.Worksheets.Item("Data")
$dataWs.Rows.Item(1).RowHeight = 31
